# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the display order / labels of "Groenlandia" (row 210) and
# "Islas Malvinas" (row 211) - they carry identical statistics, only the
# country label order in the source data changed.
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# Refresh the "last updated" timestamp string.
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 05:32"

# Kazajistan (row 30): updated case counts.
$ws.Range("B30").Value = 81720
$ws.Range("C30").Value = 1494
$ws.Range("E30").Value = 29875

# Belgica (row 37): updated case counts.
$ws.Range("B37").Value = 65727
$ws.Range("C37").Value = 528
$ws.Range("D37").Value = 17425
$ws.Range("E37").Value = 38481
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 9821

# Honduras (row 52): updated case counts.
$ws.Range("B52").Value = 38438
$ws.Range("C52").Value = 879
$ws.Range("D52").Value = 4713
$ws.Range("E52").Value = 32627
$ws.Range("G52").Value = 37
$ws.Range("H52").Value = 1098

# Haiti (row 91): updated case counts.
$ws.Range("B91").Value = 7297
$ws.Range("C91").Value = 37
$ws.Range("D91").Value = 4365
$ws.Range("E91").Value = 2775
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 157
